$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results_1")

# Row 5 - set to OOR (text) for C:L
$ws.Range("C5:L5").Value = "OOR"

# Row 6 - set to OOR (text) for C:L
$ws.Range("C6:L6").Value = "OOR"

# Row 7 - set to OOR (text) for C:L
$ws.Range("C7:L7").Value = "OOR"

# Row 8 - new numeric values
$ws.Range("C8").Value = 8019.703352289488
$ws.Range("D8").Value = 6706.350945902536
$ws.Range("E8").Value = 7139.905114274921
$ws.Range("F8").Value = 7586.149183917103
$ws.Range("G8").Value = 1.923833770379531
$ws.Range("H8").Value = 3652.252947800911
$ws.Range("I8").Value = 763.7567598375965
$ws.Range("J8").Value = 144.6126633325211
$ws.Range("K8").Value = 186.5068603760058
$ws.Range("L8").Value = 5521.252507376894

# Row 9 - new numeric values
$ws.Range("C9").Value = 7775.264560364717
$ws.Range("D9").Value = 6950.789737827307
$ws.Range("E9").Value = 7139.905114274921
$ws.Range("F9").Value = 7586.149183917103
$ws.Range("G9").Value = 1.297037211698649
$ws.Range("H9").Value = 2711.163598890543
$ws.Range("I9").Value = 763.7567598375965
$ws.Range("J9").Value = 107.3497904907753
$ws.Range("K9").Value = 182.522673224482
$ws.Range("L9").Value = 4188.036678823943

# Row 10 - new numeric values
$ws.Range("C10").Value = 7601.581740812536
$ws.Range("D10").Value = 7124.472557379489
$ws.Range("E10").Value = 7139.905114274921
$ws.Range("F10").Value = 7586.149183917103
$ws.Range("G10").Value = 0.8516750302704086
$ws.Range("H10").Value = 2042.484743614645
$ws.Range("I10").Value = 763.7567598375965
$ws.Range("J10").Value = 80.87313852891884
$ws.Range("K10").Value = 172.5285479948262
$ws.Range("L10").Value = 3337.870063514815

# Input parameter changes (B column)
$ws.Range("B14").Value = 2.31
$ws.Range("B15").Value = 0.19
$ws.Range("B16").Value = 0.018
$ws.Range("B17").Value = 0.6
$ws.Range("B18").Value = 96
$ws.Range("B20").Value = 14729
$ws.Range("B22").Value = 3.46
$ws.Range("B23").Value = 1.12
$ws.Range("B24").Value = 0.4
$ws.Range("B25").Value = 1.08
